$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "Hoofdscherm project verder uitwerken, activiteit scherm aanmaken en opmaken, instellingen scherm aanmaken, navigatie tussen deze schermen"

$ws.Range("A15").Value = "Week 4"
$ws.Range("A15").Font.Bold = $true

$ws.Range("B10").Value = "10 uur 30 minuten"
$ws.Range("B10").Font.Bold = $true

$ws.Range("B16").Value = "2 uur 15 minuten"

$ws.Range("A16").Value = 43530
$ws.Range("A16").NumberFormat = "d-mmm"

$ws.Range("B20").Select()
